# "correção nos dados e inicio da analise PNAD 2009"
#
# The original row 6 ("grandes regiões e unidades da federação") was a
# header-only row with no data; the real "norte" data row (old row 7)
# is removed and row 6 becomes the "norte" row, carrying the data that
# used to live in row 7. Every row below shifts up by one data row, and
# the last row (old row 38, "distrito federal") drops off the bottom of
# the table entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old "norte" row (row 7). This removes its text label and data
# and pulls every row below it up by one, leaving row 6 ("grandes regiões e
# unidades da federação") as the new home for what was row 7's position,
# and dropping the former last row (row 38, "distrito federal") off the
# bottom of the used range.
$ws.Rows("7:7").Delete()

# Row 6 keeps its original label slot, but it now represents "norte".
$ws.Range("A6").Value = "norte"

# Row 6 gets the data values that used to belong to the "norte" row.
$ws.Range("B6").Value = 1.53
$ws.Range("C6").Value = 1.95
$ws.Range("D6").Value = 2.99
$ws.Range("E6").Value = 2.47
$ws.Range("F6").Value = 2.12
$ws.Range("G6").Value = 1.41
